# ajustes fin de mes enero
# - Remove the "33-88 / Leidy" row (row 12) from every monthly sheet
#   (ene2025..jun2025), which shifts all subsequent rows up by one.
# - On ene2025, mark the second installment (column D, "pago2") as paid
#   (65000) for every resident, and also record the first installment
#   (column C, "pago1") for the two residents who hadn't paid it yet.
# - Update the remembered cell selection on each sheet to match where the
#   user ended up after the edit.
# - Fix a typo in the "pagos" sheet header/footer font style name.

$wb = $excel.ActiveWorkbook

# ---- ene2025 : delete the Leidy row, then true-up the payment columns ----
$wsEne = $wb.Worksheets.Item("ene2025")
$wsEne.Rows.Item(12).Delete()

# Second installment ("pago2") paid in full by everybody this month.
$wsEne.Range("D2:D22").Value = 65000
# Restore Salomon's reduced fee (row 12 after the deletion), which must
# stay at 56000 instead of the blanket 65000 used above.
$wsEne.Range("C12").Value = 56000
$wsEne.Range("D12").Value = 56000
# First installment ("pago1") settled late by Fernando (row 6) and
# Miguel (row 21, after the shift).
$wsEne.Range("C6").Value = 65000
$wsEne.Range("C21").Value = 65000

$wsEne.Range("E28").Select()

# ---- feb2025 .. jun2025 : same row removal, no payment changes ----
$wsFeb = $wb.Worksheets.Item("feb2025")
$wsFeb.Rows.Item(12).Delete()
$wsFeb.Range("A12").Select()

$wsMar = $wb.Worksheets.Item("mar2025")
$wsMar.Rows.Item(12).Delete()
$wsMar.Range("A12").Select()

$wsAbr = $wb.Worksheets.Item("abr2025")
$wsAbr.Rows.Item(12).Delete()
$wsAbr.Range("A12").Select()

$wsMay = $wb.Worksheets.Item("may2025")
$wsMay.Rows.Item(12).Delete()
$wsMay.Range("A12").Select()

$wsJun = $wb.Worksheets.Item("jun2025")
$wsJun.Rows.Item(12).Delete()
$wsJun.Range("F25").Select()

# ---- pagos : correct header/footer font style label ----
$wsPagos = $wb.Worksheets.Item("pagos")
$wsPagos.PageSetup.OddHeader = [char]38 + 'C' + [char]38 + '"Times New Roman,Normal"' + [char]38 + '12' + [char]38 + 'A'
$wsPagos.PageSetup.OddFooter = [char]38 + 'C' + [char]38 + '"Times New Roman,Normal"' + [char]38 + '12Página ' + [char]38 + 'P'

# Re-select ene2025 as the active sheet (it was the selected tab originally).
$wsEne.Activate()
